# Updates cryptos list figures (Price / Volume(1h) columns) to match
# the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.068.67"
$ws.Range("E2").Value = "  +0.37%  "

$ws.Range("D3").Value = "1.888.21"
$ws.Range("E3").Value = "  -1.32%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.68"
$ws.Range("E5").Value = "  -2.47%  "

$ws.Range("E6").Value = "  -0.14%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4597"
$ws.Range("E7").Value = "  -3.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4070"
$ws.Range("E8").Value = "  +0.33%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.72"
$ws.Range("E9").Value = "  -0.93%  "

$ws.Range("E10").Value = "  -2.23%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9925"
$ws.Range("E11").Value = "  -3.69%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.69"
$ws.Range("E12").Value = "  -3.43%  "

$ws.Range("D13").Value = "1.903.15"
$ws.Range("E13").Value = "  +0.04%  "

$ws.Range("E14").Value = "  -3.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.079"
$ws.Range("E15").Value = "  -4.18%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.44"
$ws.Range("E17").Value = "  -3.51%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001031"
$ws.Range("E18").Value = "  -2.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06555"
$ws.Range("E19").Value = "  -1.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.48"
$ws.Range("E20").Value = "  -2.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9990"
$ws.Range("E21").Value = "  -0.11%  "

$ws.Range("D22").Value = "29.089.55"
$ws.Range("E22").Value = "  +0.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.429"
$ws.Range("E23").Value = "  -2.58%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.49"
$ws.Range("E24").Value = "  +2.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.205"
$ws.Range("E25").Value = "  -2.99%  "

$ws.Range("D26").Value = "2.109.97"
$ws.Range("E26").Value = "  -0.93%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.67"
$ws.Range("E27").Value = "  -2.61%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.60"
$ws.Range("E28").Value = "  -2.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.103"
$ws.Range("E29").Value = "  -3.77%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.467"
$ws.Range("E30").Value = "  -1.35%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.78"
$ws.Range("E31").Value = "  -2.72%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.000"
$ws.Range("E32").Value = "  -1.27%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09329"
$ws.Range("E33").Value = "  -2.64%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.600"
$ws.Range("E34").Value = "  -1.55%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.410"
$ws.Range("E35").Value = "  -1.44%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.283"
$ws.Range("E36").Value = "  -2.69%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06065"
$ws.Range("E37").Value = "  -2.32%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02221"
$ws.Range("E38").Value = "  -2.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.275"
$ws.Range("E39").Value = "  -4.58%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.173"
$ws.Range("E40").Value = "  -2.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9991"
$ws.Range("E41").Value = "  -0.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5786"
$ws.Range("E42").Value = "  -4.30%  "

$ws.Range("E43").Value = "  -4.12%  "

$ws.Range("E44").Value = "  -4.41%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.258"
$ws.Range("E45").Value = "  -2.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07487"
$ws.Range("E46").Value = "  +2.45%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.05"
$ws.Range("E47").Value = "  -2.37%  "

$ws.Range("E48").Value = "  +5.18%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5454"
$ws.Range("E49").Value = "  -3.27%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.902"
$ws.Range("E50").Value = "  -4.22%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.12"
$ws.Range("E51").Value = "  -1.95%  "
